# Swap the contents of rows 11 and 12 (columns B, F:I, K:W, Z:AC) to reflect
# the re-ordering of the two Canada Premier League fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- id (column B) ---
$ws.Range("B11").Value = 6240280
$ws.Range("B12").Value = 6227815

# --- HomeTeam / AwayTeam (columns F, G) ---
$ws.Range("F11").Value = "Atletico Ottawa"
$ws.Range("G11").Value = "Vancouver FC"
$ws.Range("F12").Value = "HFX Wanderers"
$ws.Range("G12").Value = "Cavalry FC"

# --- FTHG / FTAG (columns H, I) ---
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 1

# --- odds / Asian handicap columns K:W ---
$ws.Range("K11").Value = 1.571
$ws.Range("L11").Value = 3.4
$ws.Range("M11").Value = 5.5
$ws.Range("N11").Value = 1.444
$ws.Range("O11").Value = 3.8
$ws.Range("P11").Value = 6
$ws.Range("Q11").Value = -1.25
$ws.Range("R11").Value = 1.95
$ws.Range("S11").Value = 1.85
$ws.Range("T11").Value = 2.75
$ws.Range("U11").Value = 1.975
$ws.Range("V11").Value = 1.825
$ws.Range("W11").Value = 0.444

$ws.Range("K12").Value = 2.6
$ws.Range("L12").Value = 3.2
$ws.Range("M12").Value = 2.4
$ws.Range("N12").Value = 3.3
$ws.Range("O12").Value = 3
$ws.Range("P12").Value = 2.15
$ws.Range("Q12").Value = 0.25
$ws.Range("R12").Value = 1.925
$ws.Range("S12").Value = 1.875
$ws.Range("T12").Value = 2.25
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 1.8
$ws.Range("W12").Value = 2.3

# --- profit/loss columns Z:AC (X, Y unchanged) ---
$ws.Range("Z11").Value = -0.5
$ws.Range("AA11").Value = 0.425
$ws.Range("AB11").Value = -1
$ws.Range("AC11").Value = 0.825

$ws.Range("Z12").Value = 0.925
$ws.Range("AA12").Value = -1
$ws.Range("AB12").Value = 1
$ws.Range("AC12").Value = -1
